# machineindex -> machineno 기준으로 변경 및 machineindex 컬럼 삭제.
# aging 기계도 machineindex: -1 -> machineno: 'aging'으로 변경
#
# In this workbook the "machineindex -> machineno" rename already happened
# (column A is "machineno"). What remains is to reposition the
# "machinename" column: it currently sits at column L (right after
# "duedate") and needs to move to column B (right after "machineno"),
# pushing work_start_time..duedate one column to the right. The two
# "AGING" rows (4 and 5), which previously had an empty machineindex-style
# column A, now get the literal text "AGING" in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("호기_정보")

# Insert a new, empty column at B. Everything from B onward (including the
# "machinename" column, which was L and becomes M) shifts one column right.
$ws.Columns.Item(2).Insert()

# "machinename" is now in column M (13), in rows 1-5. Cut that bounded range
# and paste it into the new column B (rows 1-5), which both removes it from
# M and populates B in one step. (Using a bounded range rather than
# Columns.Item(13).Cut(...) avoids the whole-column cut/paste needlessly
# stamping formatting across all 1,048,576 rows.)
$ws.Range("M1:M5").Cut($ws.Range("B1:B5"))

# Column M is now empty (the insert point of the cut) and superfluous;
# remove it so everything to its right shifts back left by one.
$ws.Columns.Item(13).Delete()

# Mark the two aging rows with machineno "AGING" in column A.
$ws.Range("A4").Value = "AGING"
$ws.Range("A5").Value = "AGING"
